$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B25").Value = "ArrestWarrant"
$ws.Range("C25").Value = "file "
$ws.Range("D25").Value = "true "
$ws.Range("E25").Value = "arrestWarrant"
$ws.Range("F25").Value = 50
$ws.Range("G25").Value = "P3D"
$ws.Range("H25").Value = "Review Arrest Warrant"
$ws.Range("I25").Value = "ann-acm@armedia.com,ian-acm@armedia.com,samuel-acm@armedia.com"
$ws.Range("K25").Value = "true "

$ws.Hyperlinks.Add($ws.Range("I25"), "mailto:ann-acm@armedia.com,ian-acm@armedia.com,samuel-acm@armedia.com")

$ws.Activate() | Out-Null
$ws.Range("J25").Select() | Out-Null
